$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.001.11'
$ws.Range("E2").Value = '  -0.60%  '

$ws.Range("D3").Value = '1.620.61'
$ws.Range("E3").Value = '  -0.99%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.35'
$ws.Range("E5").Value = '  -1.56%  '

$ws.Range("E6").Value = '  -0.85%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("E8").Value = '  +0.35%  '

$ws.Range("E9").Value = '  -1.30%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.91'
$ws.Range("E10").Value = '  -0.56%  '

$ws.Range("E11").Value = '  -1.20%  '

$ws.Range("D12").Value = '1.847.82'
$ws.Range("E12").Value = '  -0.99%  '

$ws.Range("D13").Value = '1.622.64'
$ws.Range("E13").Value = '  -0.95%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.10'
$ws.Range("E14").Value = '  -0.61%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.536'
$ws.Range("E15").Value = '  -0.85%  '

$ws.Range("D16").Value = '26.999.16'
$ws.Range("E16").Value = '  -0.60%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.31'
$ws.Range("E17").Value = '  -3.23%  '

$ws.Range("E18").Value = '  -0.39%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '213.44'
$ws.Range("E19").Value = '  -1.41%  '

$ws.Range("E20").Value = '  +0.01%  '

$ws.Range("E21").Value = '  -0.16%  '

$ws.Range("E22").Value = '  -2.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.36'
$ws.Range("E23").Value = '  -7.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.93'
$ws.Range("E24").Value = '  -1.97%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.09'
$ws.Range("E25").Value = '  -0.29%  '

$ws.Range("E26").Value = '  +1.17%  '

$ws.Range("E27").Value = '  +0.02%  '

$ws.Range("E28").Value = '  -3.70%  '

$ws.Range("E29").Value = '  -1.04%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0507'
$ws.Range("E30").Value = '  +0.06%  '

$ws.Range("E31").Value = '  -1.24%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.29'
$ws.Range("E32").Value = '  -2.65%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.705'
$ws.Range("E33").Value = '  +28.48%  '

$ws.Range("E34").Value = '  -1.09%  '

$ws.Range("D35").Value = '1.333.74'
$ws.Range("E35").Value = '  +2.44%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.55'
$ws.Range("E36").Value = '  -0.71%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.44'
$ws.Range("E37").Value = '  -0.52%  '

$ws.Range("E38").Value = '  -0.69%  '

$ws.Range("E39").Value = '  -1.70%  '

$ws.Range("E40").Value = '  +0.01%  '

$ws.Range("E41").Value = '  -1.33%  '

$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.35'
$ws.Range("E42").Value = '  +0.22%  '

$ws.Range("B43").Value = 'MXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.18'
$ws.Range("E43").Value = '  -2.62%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.80'
$ws.Range("E44").Value = '  +2.70%  '

$ws.Range("D45").Value = '1.759.09'
$ws.Range("E45").Value = '  -0.97%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.87'
$ws.Range("E46").Value = '  -1.48%  '

$ws.Range("E47").Value = '  +2.04%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.847'
$ws.Range("E48").Value = '  +16.18%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0517'
$ws.Range("E49").Value = '  +0.38%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0989'
$ws.Range("E50").Value = '  +3.42%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.53'
$ws.Range("E51").Value = '  -1.52%  '

# Reset number-formatted cells back to default (General) style, matching original formatting,
# by pasting the format from an adjacent default-styled cell in the same row.
$ws.Range("B5").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("B10").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("B14").Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("B15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("B17").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("B19").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Range("B23").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("B24").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Range("B25").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("B30").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("B32").Copy()
$ws.Range("D32").PasteSpecial(-4122)
$ws.Range("B33").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("B36").Copy()
$ws.Range("D36").PasteSpecial(-4122)
$ws.Range("B37").Copy()
$ws.Range("D37").PasteSpecial(-4122)
$ws.Range("B42").Copy()
$ws.Range("D42").PasteSpecial(-4122)
$ws.Range("B43").Copy()
$ws.Range("D43").PasteSpecial(-4122)
$ws.Range("B44").Copy()
$ws.Range("D44").PasteSpecial(-4122)
$ws.Range("B46").Copy()
$ws.Range("D46").PasteSpecial(-4122)
$ws.Range("B48").Copy()
$ws.Range("D48").PasteSpecial(-4122)
$ws.Range("B49").Copy()
$ws.Range("D49").PasteSpecial(-4122)
$ws.Range("B50").Copy()
$ws.Range("D50").PasteSpecial(-4122)
$ws.Range("B51").Copy()
$ws.Range("D51").PasteSpecial(-4122)
$excel.CutCopyMode = 0